$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 28.878682
$ws.Cells.Item(2, 8).Value = 86.63604599999999
$ws.Cells.Item(2, 9).Value = 0.2978181586389064
$ws.Cells.Item(2, 10).Value = 0.2978181586389064
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 8.131233999999999
$ws.Cells.Item(2, 14).Value = 24.393702
$ws.Cells.Item(2, 15).Value = 0.02090995573015822
$ws.Cells.Item(2, 16).Value = 0.02090995573015823
$ws.Cells.Item(2, 17).Value = 234.8193209535879
$ws.Cells.Item(2, 18).Value = 2113.373888582292
$ws.Cells.Item(2, 19).Value = 0.006227364512776771
$ws.Cells.Item(2, 20).Value = 0.006227364512776772

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 28.878682
$ws.Cells.Item(3, 8).Value = 86.63604599999999
$ws.Cells.Item(3, 9).Value = 0.2978181586389064
$ws.Cells.Item(3, 10).Value = 0.2978181586389064
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 243.3763986666667
$ws.Cells.Item(3, 14).Value = 730.1291960000001
$ws.Cells.Item(3, 15).Value = 0.625857000534647
$ws.Cells.Item(3, 16).Value = 0.6258570005346471
$ws.Cells.Item(3, 17).Value = 7028.38962339989
$ws.Cells.Item(3, 18).Value = 63255.50661059902
$ws.Cells.Item(3, 19).Value = 0.1863915794704976
$ws.Cells.Item(3, 20).Value = 0.1863915794704976

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 28.878682
$ws.Cells.Item(4, 8).Value = 86.63604599999999
$ws.Cells.Item(4, 9).Value = 0.2978181586389064
$ws.Cells.Item(4, 10).Value = 0.2978181586389064
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 103.9426383333333
$ws.Cells.Item(4, 14).Value = 311.827915
$ws.Cells.Item(4, 15).Value = 0.2672947262403034
$ws.Cells.Item(4, 16).Value = 0.2672947262403035
$ws.Cells.Item(4, 17).Value = 3001.726398669343
$ws.Cells.Item(4, 18).Value = 27015.53758802409
$ws.Cells.Item(4, 19).Value = 0.07960522318277774
$ws.Cells.Item(4, 20).Value = 0.07960522318277775

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 28.878682
$ws.Cells.Item(5, 8).Value = 86.63604599999999
$ws.Cells.Item(5, 9).Value = 0.2978181586389064
$ws.Cells.Item(5, 10).Value = 0.2978181586389064
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 33.41874933333333
$ws.Cells.Item(5, 14).Value = 100.256248
$ws.Cells.Item(5, 15).Value = 0.08593831749489127
$ws.Cells.Item(5, 16).Value = 0.08593831749489128
$ws.Cells.Item(5, 17).Value = 965.0894348350452
$ws.Cells.Item(5, 18).Value = 8685.804913515407
$ws.Cells.Item(5, 19).Value = 0.02559399147285423
$ws.Cells.Item(5, 20).Value = 0.02559399147285423

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 45.41653666666667
$ws.Cells.Item(6, 8).Value = 136.24961
$ws.Cells.Item(6, 9).Value = 0.4683686506822938
$ws.Cells.Item(6, 10).Value = 0.4683686506822937
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 8.131233999999999
$ws.Cells.Item(6, 14).Value = 24.393702
$ws.Cells.Item(6, 15).Value = 0.02090995573015822
$ws.Cells.Item(6, 16).Value = 0.02090995573015823
$ws.Cells.Item(6, 17).Value = 369.2924871062467
$ws.Cells.Item(6, 18).Value = 3323.63238395622
$ws.Cells.Item(6, 19).Value = 0.009793567751160704
$ws.Cells.Item(6, 20).Value = 0.009793567751160704

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 45.41653666666667
$ws.Cells.Item(7, 8).Value = 136.24961
$ws.Cells.Item(7, 9).Value = 0.4683686506822938
$ws.Cells.Item(7, 10).Value = 0.4683686506822937
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 243.3763986666667
$ws.Cells.Item(7, 14).Value = 730.1291960000001
$ws.Cells.Item(7, 15).Value = 0.625857000534647
$ws.Cells.Item(7, 16).Value = 0.6258570005346471
$ws.Cells.Item(7, 17).Value = 11053.31313384595
$ws.Cells.Item(7, 18).Value = 99479.81820461359
$ws.Cells.Item(7, 19).Value = 0.2931317988604802
$ws.Cells.Item(7, 20).Value = 0.2931317988604802

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 45.41653666666667
$ws.Cells.Item(8, 8).Value = 136.24961
$ws.Cells.Item(8, 9).Value = 0.4683686506822938
$ws.Cells.Item(8, 10).Value = 0.4683686506822937
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 103.9426383333333
$ws.Cells.Item(8, 14).Value = 311.827915
$ws.Cells.Item(8, 15).Value = 0.2672947262403034
$ws.Cells.Item(8, 16).Value = 0.2672947262403035
$ws.Cells.Item(8, 17).Value = 4720.714645095906
$ws.Cells.Item(8, 18).Value = 42486.43180586316
$ws.Cells.Item(8, 19).Value = 0.125192470263664
$ws.Cells.Item(8, 20).Value = 0.125192470263664

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 45.41653666666667
$ws.Cells.Item(9, 8).Value = 136.24961
$ws.Cells.Item(9, 9).Value = 0.4683686506822938
$ws.Cells.Item(9, 10).Value = 0.4683686506822937
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 33.41874933333333
$ws.Cells.Item(9, 14).Value = 100.256248
$ws.Cells.Item(9, 15).Value = 0.08593831749489127
$ws.Cells.Item(9, 16).Value = 0.08593831749489128
$ws.Cells.Item(9, 17).Value = 1517.763854451476
$ws.Cells.Item(9, 18).Value = 13659.87469006328
$ws.Cells.Item(9, 19).Value = 0.04025081380698878
$ws.Cells.Item(9, 20).Value = 0.04025081380698878

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 6.995916999999999
$ws.Cells.Item(10, 8).Value = 20.987751
$ws.Cells.Item(10, 9).Value = 0.0721470293876508
$ws.Cells.Item(10, 10).Value = 0.07214702938765079
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 8.131233999999999
$ws.Cells.Item(10, 14).Value = 24.393702
$ws.Cells.Item(10, 15).Value = 0.02090995573015822
$ws.Cells.Item(10, 16).Value = 0.02090995573015823
$ws.Cells.Item(10, 17).Value = 56.88543817157799
$ws.Cells.Item(10, 18).Value = 511.9689435442019
$ws.Cells.Item(10, 19).Value = 0.001508591190558203
$ws.Cells.Item(10, 20).Value = 0.001508591190558203

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 6.995916999999999
$ws.Cells.Item(11, 8).Value = 20.987751
$ws.Cells.Item(11, 9).Value = 0.0721470293876508
$ws.Cells.Item(11, 10).Value = 0.07214702938765079
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 243.3763986666667
$ws.Cells.Item(11, 14).Value = 730.1291960000001
$ws.Cells.Item(11, 15).Value = 0.625857000534647
$ws.Cells.Item(11, 16).Value = 0.6258570005346471
$ws.Cells.Item(11, 17).Value = 1702.641084830911
$ws.Cells.Item(11, 18).Value = 15323.7697634782
$ws.Cells.Item(11, 19).Value = 0.04515372341004016
$ws.Cells.Item(11, 20).Value = 0.04515372341004016

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 6.995916999999999
$ws.Cells.Item(12, 8).Value = 20.987751
$ws.Cells.Item(12, 9).Value = 0.0721470293876508
$ws.Cells.Item(12, 10).Value = 0.07214702938765079
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 103.9426383333333
$ws.Cells.Item(12, 14).Value = 311.827915
$ws.Cells.Item(12, 15).Value = 0.2672947262403034
$ws.Cells.Item(12, 16).Value = 0.2672947262403035
$ws.Cells.Item(12, 17).Value = 727.1740705410183
$ws.Cells.Item(12, 18).Value = 6544.566634869165
$ws.Cells.Item(12, 19).Value = 0.01928452046922325
$ws.Cells.Item(12, 20).Value = 0.01928452046922325

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 6.995916999999999
$ws.Cells.Item(13, 8).Value = 20.987751
$ws.Cells.Item(13, 9).Value = 0.0721470293876508
$ws.Cells.Item(13, 10).Value = 0.07214702938765079
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 33.41874933333333
$ws.Cells.Item(13, 14).Value = 100.256248
$ws.Cells.Item(13, 15).Value = 0.08593831749489127
$ws.Cells.Item(13, 16).Value = 0.08593831749489128
$ws.Cells.Item(13, 17).Value = 233.7947965798053
$ws.Cells.Item(13, 18).Value = 2104.153169218248
$ws.Cells.Item(13, 19).Value = 0.006200194317829185
$ws.Cells.Item(13, 20).Value = 0.006200194317829185

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 15.67636333333333
$ws.Cells.Item(14, 8).Value = 47.02909
$ws.Cells.Item(14, 9).Value = 0.161666161291149
$ws.Cells.Item(14, 10).Value = 0.161666161291149
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 8.131233999999999
$ws.Cells.Item(14, 14).Value = 24.393702
$ws.Cells.Item(14, 15).Value = 0.02090995573015822
$ws.Cells.Item(14, 16).Value = 0.02090995573015823
$ws.Cells.Item(14, 17).Value = 127.4681785323533
$ws.Cells.Item(14, 18).Value = 1147.21360679118
$ws.Cells.Item(14, 19).Value = 0.003380432275662545
$ws.Cells.Item(14, 20).Value = 0.003380432275662545

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 15.67636333333333
$ws.Cells.Item(15, 8).Value = 47.02909
$ws.Cells.Item(15, 9).Value = 0.161666161291149
$ws.Cells.Item(15, 10).Value = 0.161666161291149
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 243.3763986666667
$ws.Cells.Item(15, 14).Value = 730.1291960000001
$ws.Cells.Item(15, 15).Value = 0.625857000534647
$ws.Cells.Item(15, 16).Value = 0.6258570005346471
$ws.Cells.Item(15, 17).Value = 3815.256852256849
$ws.Cells.Item(15, 18).Value = 34337.31167031164
$ws.Cells.Item(15, 19).Value = 0.101179898793629
$ws.Cells.Item(15, 20).Value = 0.101179898793629

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 15.67636333333333
$ws.Cells.Item(16, 8).Value = 47.02909
$ws.Cells.Item(16, 9).Value = 0.161666161291149
$ws.Cells.Item(16, 10).Value = 0.161666161291149
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 103.9426383333333
$ws.Cells.Item(16, 14).Value = 311.827915
$ws.Cells.Item(16, 15).Value = 0.2672947262403034
$ws.Cells.Item(16, 16).Value = 0.2672947262403035
$ws.Cells.Item(16, 17).Value = 1629.442564338594
$ws.Cells.Item(16, 18).Value = 14664.98307904735
$ws.Cells.Item(16, 19).Value = 0.04321251232463841
$ws.Cells.Item(16, 20).Value = 0.04321251232463842

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 15.67636333333333
$ws.Cells.Item(17, 8).Value = 47.02909
$ws.Cells.Item(17, 9).Value = 0.161666161291149
$ws.Cells.Item(17, 10).Value = 0.161666161291149
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 33.41874933333333
$ws.Cells.Item(17, 14).Value = 100.256248
$ws.Cells.Item(17, 15).Value = 0.08593831749489127
$ws.Cells.Item(17, 16).Value = 0.08593831749489128
$ws.Cells.Item(17, 17).Value = 523.8844566949244
$ws.Cells.Item(17, 18).Value = 4714.960110254319
$ws.Cells.Item(17, 19).Value = 0.01389331789721907
$ws.Cells.Item(17, 20).Value = 0.01389331789721907
